# LOM3018.docx edit: rotate the body text blocks between sections while
# keeping all headings / bold labels / paragraph formatting untouched.
#
# Every "value" text block (Objetivos, Docente names, Programa resumido,
# Programa, the Metodo/Criterio/Norma values, and the Bibliografia
# entries) moves to a different paragraph/run slot - a full rotation
# where every source text is also somebody else's destination. To do
# this safely and deterministically:
#   1) each of the 9 source locations is stamped with a unique
#      placeholder token (one Find/Replace per location). Locations that
#      hold more than one <w:t> (joined by manual line breaks, i.e.
#      <w:br/>) are matched in a single Find call that spans the whole
#      run group (using [char]11 for the break), so the whole group
#      collapses to one placeholder with no leftover <w:br/>.
#   2) each placeholder is then replaced with the final text for that
#      slot, again using [char]11 to (re)introduce manual line breaks
#      exactly where the target needs them.
# Because every search token in phase 2 is a unique placeholder, there
# is no risk of a replacement being applied at the wrong spot.

$d = $word.ActiveDocument

function Replace-InDoc($oldText, $newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find.Execute failed to locate: $oldText"
    }
}

$LF = [char]11

# ---- Phase 1: stamp a unique placeholder over each of the 9 source slots ----

# para 6 - Objetivos (single run)
Replace-InDoc "A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar." "@@SLOT_OBJ@@"

# para 8 - Docente(s), two runs joined by a manual break -> one placeholder
$oldDocente = "984972 - Hugo Ricardo Zschommler Sandim" + $LF + "7459752 - Maria Ismenia Sodero Toledo Faria"
Replace-InDoc $oldDocente "@@SLOT_DOCENTE@@"

# para 10 - Programa resumido, two runs joined by a manual break -> one placeholder
$oldResumido = "1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia." + $LF + "Em todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos"
Replace-InDoc $oldResumido "@@SLOT_RESUMIDO@@"

# para 12 - Programa (single run)
Replace-InDoc "As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos. Viagem didática complementar." "@@SLOT_PROGRAMA@@"

# para 14 - Avaliação values (Método / Critério / Norma de recuperação)
Replace-InDoc "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras" "@@SLOT_METODO@@"

$oldCriterioVal = "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
Replace-InDoc $oldCriterioVal "@@SLOT_CRITERIO@@"

Replace-InDoc "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação." "@@SLOT_NORMA@@"

# para 16 - Bibliografia, five <w:t> joined by manual breaks inside one run -> one placeholder
$oldBibliografia = "1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010." + $LF + `
    "2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006." + $LF + `
    "2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. " + $LF + `
    "4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985." + $LF + `
    "5) Artigos científicos"
Replace-InDoc $oldBibliografia "@@SLOT_BIBLIOGRAFIA@@"

# ---- Phase 2: write the final content into each slot (by its placeholder) ----

# para 6 now holds the old "Programa resumido" text (two parts joined by
# a manual line break).
$newObjSlot = "1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia." + $LF + "Em todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos"
Replace-InDoc "@@SLOT_OBJ@@" $newObjSlot

# para 8 now holds the old "Objetivos" text (run 1) and the old
# "Programa" text (run 2), joined by a manual line break.
$newDocenteSlot = "A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar." + $LF + "As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos. Viagem didática complementar."
Replace-InDoc "@@SLOT_DOCENTE@@" $newDocenteSlot

# para 10 now holds the old "Método" text as a single run (no break).
Replace-InDoc "@@SLOT_RESUMIDO@@" "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"

# para 12 now holds the old "Critério" text.
Replace-InDoc "@@SLOT_PROGRAMA@@" "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# para 14, "Método:" value now holds the old "Norma de recuperação" text.
Replace-InDoc "@@SLOT_METODO@@" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# para 14, "Critério:" value now holds the old Bibliografia entries
# (joined by manual line breaks, same structure as before).
$newCriterioSlot = "1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010." + $LF + `
    "2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006." + $LF + `
    "2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. " + $LF + `
    "4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985." + $LF + `
    "5) Artigos científicos"
Replace-InDoc "@@SLOT_CRITERIO@@" $newCriterioSlot

# para 14, "Norma de recuperação:" value now holds the old first Docente
# text.
Replace-InDoc "@@SLOT_NORMA@@" "984972 - Hugo Ricardo Zschommler Sandim"

# para 16 now holds the old second Docente text as a single run.
Replace-InDoc "@@SLOT_BIBLIOGRAFIA@@" "7459752 - Maria Ismenia Sodero Toledo Faria"

Write-Output "Done."
